$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item("TextBox 34")
$tr = $shp.TextFrame.TextRange
Write-Host $tr.Text
Write-Host $tr.Runs().Count
for ($i = 1; $i -le $tr.Runs().Count; $i++) {
    $r = $tr.Runs($i)
    Write-Host "Run $i : [$($r.Text)]"
}
